# "Generate Report for handback" - append the handback row for
# 286d61d8-6c7c-42a0-8284-698db6bcb61a to the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$uid  = "286d61d8-6c7c-42a0-8284-698db6bcb61a"
$hash = "d852b6296bfaf02417579a183a57733436a31a83"

$mdName = "$uid.md"
$statusInSync = "Handed back: in sync with en-US"
$include = "Include"

# ---------------------------------------------------------------------------
# Overview sheet - row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/99636c06de6255875b54fcb7bd6cc5e695c3f96f/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------------
# Per-language detail sheets (zh-cn / de-de) - row 4
# ---------------------------------------------------------------------------
$langs = "zh-cn", "de-de"

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang)

    $xlfName = "$uid.$hash.$lang.xlf"

    # A4: Source File Name -> md in the main oltest repo
    $ws.Hyperlinks.Add(
        $ws.Range("A4"),
        "https://github.com/OpenLocalizationTest/oltest/blob/99636c06de6255875b54fcb7bd6cc5e695c3f96f/e2e/$mdName",
        "",
        "",
        $mdName
    )

    # B4: Status
    $ws.Range("B4").Value = $statusInSync

    # C4: Correspond Handoff File -> handoff xlf
    $ws.Hyperlinks.Add(
        $ws.Range("C4"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e964cb77306ae0c76c02278cb2f20adc16505f60/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/xinjiang/ht/$xlfName",
        "",
        "",
        $xlfName
    )

    # D4: Correspond Handoff Datetime
    $ws.Range("D4").Value = "2016-02-15 08:05:04"
    $ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    # E4: Target File -> md in the localized oltest.<lang> repo
    $ws.Hyperlinks.Add(
        $ws.Range("E4"),
        "https://github.com/OpenLocalizationTestOrg/oltest.$lang/blob/dc71e8e48a9ba7fcdf82e2ef7470528c51bbfb28/e2e/$mdName",
        "",
        "",
        $mdName
    )

    # F4: Correspond Handback File -> same xlf as handoff, now in olhandback repo
    $ws.Hyperlinks.Add(
        $ws.Range("F4"),
        "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2131c719034c3ffa2504e436ee0df122f81c3cc3/ol-handback/OpenLocalizationTestOrg/oltest.$lang/xinjiang/ht/$xlfName",
        "",
        "",
        $xlfName
    )

    # G4: Correspond Handback DateTime (plain text, no date number format)
    $ws.Range("G4").Value = "2016-02-15 08:05:51"

    # H4: Handoff Reason
    $ws.Range("H4").Value = $include
}

# Fix up the per-language datetimes (zh-cn / de-de differ slightly per the diff).
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-15 08:05:04"
$wsZh.Range("G4").Value = "2016-02-15 08:05:51"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-15 08:05:18"
$wsDe.Range("G4").Value = "2016-02-15 08:06:17"
